# Insert a new slide (duplicate of the "I am bob... / But I have to do this,
# now let me pause my current working project ui here." slide, together with
# its picture) right after slide 1, i.e. at position 2.
#
# Slides 5 and 6 of the original deck already carry the exact title / body /
# picture combination the new slide needs, so the new slide is produced by
# duplicating slide 6 and moving the duplicate up to position 2 - this keeps
# all text runs, formatting and the embedded picture (image4.png) identical
# to what PowerPoint itself would produce when a user duplicates a slide and
# drags it to a new spot.

$p = $ppt.ActivePresentation

$source = $p.Slides.Item(6)
$newSlide = $source.Duplicate()
$newSlide.MoveTo(2)
